$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 361.5
$ws.Range("I6").Value = 160
$ws.Range("J6").Value = 630.1667
$ws.Range("K6").Value = 480
$ws.Range("L6").Value = 1890.5001
$ws.Range("M6").Value = -368
$ws.Range("N6").Value = -2114.5001

$ws.Range("H31").Value = 600
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 600
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 1800
$ws.Range("M31").Value = $null
$ws.Range("N31").Value = -2260

$ws.Range("H63").Value = 16265.667
$ws.Range("I63").Value = 6900
$ws.Range("J63").Value = 34997
$ws.Range("K63").Value = 6900
$ws.Range("L63").Value = 34997
$ws.Range("M63").Value = -6276
$ws.Range("N63").Value = -36245

$ws.Range("H66").Value = 16265.667
$ws.Range("I66").Value = 6900
$ws.Range("J66").Value = 34997
$ws.Range("K66").Value = 20700
$ws.Range("L66").Value = 104991
$ws.Range("M66").Value = -17580
$ws.Range("N66").Value = -111231

$ws.Range("H88").Value = 2875.0557
$ws.Range("I88").Value = 1700
$ws.Range("J88").Value = 3110.0667
$ws.Range("K88").Value = 1700
$ws.Range("L88").Value = 3110.0667
$ws.Range("M88").Value = -1294
$ws.Range("N88").Value = -3922.0667

$ws.Range("H91").Value = 2875.0557
$ws.Range("I91").Value = 1700
$ws.Range("J91").Value = 3110.0667
$ws.Range("K91").Value = 1700
$ws.Range("L91").Value = 3110.0667
$ws.Range("M91").Value = -296
$ws.Range("N91").Value = -5918.066699999999

$ws.Range("H129").Value = 3240.2195
$ws.Range("J129").Value = 1099.5
$ws.Range("L129").Value = 3298.5
$ws.Range("N129").Value = -13298.5

$ws.Range("H137").Value = 1605.1818
$ws.Range("J137").Value = 1352.75
$ws.Range("L137").Value = 4058.25
$ws.Range("N137").Value = -9158.25

$ws.Range("H138").Value = 1559.079
$ws.Range("J138").Value = 1558.8
$ws.Range("L138").Value = 4676.4
$ws.Range("N138").Value = -14956.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2310
$ws.Range("J63").Value = 3100
$ws.Range("L63").Value = 3100
$ws.Range("N63").Value = -4472

$ws.Range("H66").Value = 2310
$ws.Range("J66").Value = 3100
$ws.Range("L66").Value = 15500
$ws.Range("N66").Value = -22364

$ws.Range("H86").Value = 45485
$ws.Range("J86").Value = 45485
$ws.Range("L86").Value = 45485
$ws.Range("N86").Value = -47857

$ws.Range("H89").Value = 45485
$ws.Range("J89").Value = 45485
$ws.Range("L89").Value = 136455
$ws.Range("N89").Value = -148311

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 499.33334
$ws.Range("I22").Value = 519.6
$ws.Range("J22").Value = 398
$ws.Range("K22").Value = 519.6
$ws.Range("L22").Value = 398
$ws.Range("M22").Value = -346.6
$ws.Range("N22").Value = -744

$ws.Range("H99").Value = 1361.7
$ws.Range("J99").Value = 1582.6923
$ws.Range("L99").Value = 1582.6923
$ws.Range("N99").Value = -4578.6923

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 1000
$ws.Range("J15").Value = 1000
$ws.Range("L15").Value = 1000
$ws.Range("N15").Value = -1340

$ws.Range("H82").Value = 33135.75
$ws.Range("J82").Value = 33135.75
$ws.Range("L82").Value = 33135.75
$ws.Range("N82").Value = -33857.75

$ws.Range("H85").Value = 33135.75
$ws.Range("J85").Value = 33135.75
$ws.Range("L85").Value = 33135.75
$ws.Range("N85").Value = -35631.75

$ws.Range("H132").Value = 3121.7112
$ws.Range("I132").Value = 2782.861
$ws.Range("J132").Value = 4477.1113
$ws.Range("K132").Value = 8348.582999999999
$ws.Range("L132").Value = 13431.3339
$ws.Range("M132").Value = -5818.582999999999
$ws.Range("N132").Value = -18491.3339

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 500
$ws.Range("I16").Value = 500
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1500
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1327
$ws.Range("N16").Value = $null

$ws.Range("H24").Value = 60
$ws.Range("J24").Value = 60
$ws.Range("L24").Value = 180
$ws.Range("N24").Value = -640

$ws.Range("H131").Value = 816.3200000000001
$ws.Range("J131").Value = 816.3200000000001
$ws.Range("L131").Value = 2448.96
$ws.Range("N131").Value = -12528.96

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3725.724
$ws.Range("I7").Value = 4119.4
$ws.Range("K7").Value = 4119.4
$ws.Range("M7").Value = -4007.4

$ws.Range("H22").Value = 5433.1665
$ws.Range("J22").Value = 3250
$ws.Range("L22").Value = 3250
$ws.Range("N22").Value = -3840

$ws.Range("H27").Value = 5433.1665
$ws.Range("J27").Value = 3250
$ws.Range("L27").Value = 3250
$ws.Range("N27").Value = -3464

$ws.Range("H61").Value = 1541.8462
$ws.Range("I61").Value = 1470.4706
$ws.Range("J61").Value = 1676.6666
$ws.Range("K61").Value = 1470.4706
$ws.Range("L61").Value = 1676.6666
$ws.Range("M61").Value = -1268.4706
$ws.Range("N61").Value = -2080.6666

$ws.Range("H113").Value = 1541.8462
$ws.Range("I113").Value = 1470.4706
$ws.Range("J113").Value = 1676.6666
$ws.Range("K113").Value = 1470.4706
$ws.Range("L113").Value = 1676.6666
$ws.Range("M113").Value = 699.5293999999999
$ws.Range("N113").Value = -6016.6666

$ws.Range("H126").Value = 3725.724
$ws.Range("I126").Value = 4119.4
$ws.Range("K126").Value = 12358.2
$ws.Range("M126").Value = -9888.199999999999

$ws.Range("H132").Value = 4223.6772
$ws.Range("I132").Value = 3794.077
$ws.Range("J132").Value = 6457.6
$ws.Range("K132").Value = 11382.231
$ws.Range("L132").Value = 19372.8
$ws.Range("M132").Value = -8852.231
$ws.Range("N132").Value = -24432.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 13683.667
$ws.Range("I61").Value = 14525.5
$ws.Range("K61").Value = 14525.5
$ws.Range("M61").Value = -14233.5

$ws.Range("H95").Value = 23344
$ws.Range("J95").Value = 23344
$ws.Range("L95").Value = 23344
$ws.Range("N95").Value = -28836

$ws.Range("H132").Value = 2693.4119
$ws.Range("J132").Value = 3980
$ws.Range("L132").Value = 11940
$ws.Range("N132").Value = -17000
